$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9109108448028564,
    0.8708708882331848,
    0.9309308528900146,
    0.9409409165382385,
    0.9109108448028564,
    0.9409409165382385,
    0.9409409165382385,
    0.9509509205818176,
    0.9209208488464355,
    0.9309308528900146,
    0.9309308528900146,
    0.9309308528900146,
    0.9409409165382385,
    0.9609609246253967,
    0.9309308528900146,
    0.9609609246253967,
    0.9609609246253967,
    0.9509509205818176,
    0.9409409165382385,
    0.9309308528900146,
    0.9809809327125549,
    0.9209208488464355,
    0.9509509205818176,
    0.9409409165382385,
    0.9309308528900146,
    0.9009008407592773,
    0.9409409165382385,
    0.9509509205818176,
    0.9609609246253967,
    0.9209208488464355,
    0.9409409165382385,
    0.9409409165382385,
    0.9409409165382385,
    0.9209208488464355,
    0.9209208488464355,
    0.9309308528900146,
    0.9309308528900146,
    0.9309308528900146,
    0.9509509205818176,
    0.9209208488464355,
    0.9309308528900146,
    0.9309308528900146,
    0.9409409165382385,
    0.9309308528900146,
    0.9409409165382385,
    0.9409409165382385,
    0.9609609246253967,
    0.9709709286689758,
    0.9509509205818176,
    0.9509509205818176,
    0.9609609246253967,
    0.9509509205818176,
    0.9409409165382385
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}
